# Auto-generated Excel COM-interop script that applies the Hyperion_Profits data refresh
# (scheduled runner update) described by the commit diff: updates to computed market-board
# columns H:N (currentAveragePrice*, LevePrice*, LeveProfit*) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# row edit (@@ -2616,25 +2616,25 @@)
$ws.Range("H40").Value = 3634.4167
$ws.Range("I40").Value = 1783.8334
$ws.Range("J40").Value = 5485
$ws.Range("K40").Value = 1783.8334
$ws.Range("L40").Value = 5485
$ws.Range("M40").Value = -1608.8334
$ws.Range("N40").Value = -5835
# row edit (@@ -2772,25 +2772,25 @@)
$ws.Range("H43").Value = 1332.8462
$ws.Range("I43").Value = 1343.5
$ws.Range("J43").Value = 1297.3334
$ws.Range("K43").Value = 1343.5
$ws.Range("L43").Value = 1297.3334
$ws.Range("M43").Value = -1274.5
$ws.Range("N43").Value = -1435.3334
# row edit (@@ -3709,22 +3709,22 @@)
$ws.Range("H62").Value = 4877.857
$ws.Range("I62").Value = 3526.4546
$ws.Range("K62").Value = 3526.4546
$ws.Range("M62").Value = -2902.4546
# row edit (@@ -3859,22 +3859,22 @@)
$ws.Range("H65").Value = 4877.857
$ws.Range("I65").Value = 3526.4546
$ws.Range("K65").Value = 17632.273
$ws.Range("M65").Value = -14512.273
# row edit (@@ -6268,25 +6268,25 @@)
$ws.Range("H113").Value = 5547.8125
$ws.Range("I113").Value = 4499
$ws.Range("J113").Value = 5897.4165
$ws.Range("K113").Value = 4499
$ws.Range("L113").Value = 5897.4165
$ws.Range("M113").Value = -1245
$ws.Range("N113").Value = -12405.4165
# row edit (@@ -6853,25 +6853,25 @@)
$ws.Range("H125").Value = 13336022
$ws.Range("J125").Value = 17546896
$ws.Range("L125").Value = 157922064
$ws.Range("N125").Value = -157926984
# row edit (@@ -7205,22 +7205,22 @@)
$ws.Range("H132").Value = 29414984
$ws.Range("I132").Value = 43482030
$ws.Range("K132").Value = 130446090
$ws.Range("M132").Value = -130443560
# row edit (@@ -7453,25 +7453,25 @@)
$ws.Range("H137").Value = 65990.39
$ws.Range("I137").Value = 86631.57000000001
$ws.Range("J137").Value = 4066.8572
$ws.Range("K137").Value = 259894.71
$ws.Range("L137").Value = 12200.5716
$ws.Range("M137").Value = -257344.71
$ws.Range("N137").Value = -17300.5716
# row edit (@@ -7505,25 +7505,25 @@)
$ws.Range("H138").Value = 3323.8572
$ws.Range("I138").Value = 1600.6842
$ws.Range("J138").Value = 4747.3477
$ws.Range("K138").Value = 4802.0526
$ws.Range("L138").Value = 14242.0431
$ws.Range("M138").Value = 337.9474
$ws.Range("N138").Value = -24522.0431

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# row edit (@@ -9923,22 +9923,22 @@)
$ws.Range("H45").Value = 5141029
$ws.Range("I45").Value = 6851897
$ws.Range("K45").Value = 6851897
$ws.Range("M45").Value = -6851520
# row edit (@@ -12435,22 +12435,22 @@)
$ws.Range("H97").Value = 1798594.5
$ws.Range("I97").Value = 1798594.5
$ws.Range("K97").Value = 1798594.5
$ws.Range("M97").Value = -1798098.5
# row edit (@@ -13639,22 +13639,22 @@)
$ws.Range("H122").Value = 1605391.4
$ws.Range("I122").Value = 2878.9
$ws.Range("K122").Value = 8636.700000000001
$ws.Range("M122").Value = -6186.700000000001
# row edit (@@ -14509,19 +14509,22 @@)
$ws.Range("H140").Value = 49000
$ws.Range("J140").Value = 49000
$ws.Range("L140").Value = 49000
$ws.Range("N140").Value = -59360

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# row edit (@@ -18799,25 +18802,25 @@)
$ws.Range("H86").Value = 7154227.5
$ws.Range("I86").Value = 7704476
$ws.Range("J86").Value = 996
$ws.Range("K86").Value = 7704476
$ws.Range("L86").Value = 996
$ws.Range("M86").Value = -7703353
$ws.Range("N86").Value = -3242
# row edit (@@ -18946,25 +18949,25 @@)
$ws.Range("H89").Value = 7154227.5
$ws.Range("I89").Value = 7704476
$ws.Range("J89").Value = 996
$ws.Range("K89").Value = 38522380
$ws.Range("L89").Value = 4980
$ws.Range("M89").Value = -38516764
$ws.Range("N89").Value = -16212
# row edit (@@ -21103,25 +21106,25 @@)
$ws.Range("H134").Value = 5030.7
$ws.Range("I134").Value = 2048.077
$ws.Range("J134").Value = 10569.857
$ws.Range("K134").Value = 6144.231000000001
$ws.Range("L134").Value = 31709.571
$ws.Range("M134").Value = -3609.231000000001
$ws.Range("N134").Value = -36779.571

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# row edit (@@ -22995,25 +22998,25 @@)
$ws.Range("H31").Value = 16194.162
$ws.Range("I31").Value = 1951.7059
$ws.Range("J31").Value = 20441.912
$ws.Range("K31").Value = 1951.7059
$ws.Range("L31").Value = 20441.912
$ws.Range("M31").Value = -1656.7059
$ws.Range("N31").Value = -21031.912
# row edit (@@ -23148,25 +23151,25 @@)
$ws.Range("H34").Value = 16194.162
$ws.Range("I34").Value = 1951.7059
$ws.Range("J34").Value = 20441.912
$ws.Range("K34").Value = 1951.7059
$ws.Range("L34").Value = 20441.912
$ws.Range("M34").Value = -1749.7059
$ws.Range("N34").Value = -20845.912
# row edit (@@ -23200,22 +23203,22 @@)
$ws.Range("H35").Value = 3383.1667
$ws.Range("I35").Value = 1859.8
$ws.Range("K35").Value = 1859.8
$ws.Range("M35").Value = -1565.8
# row edit (@@ -23402,22 +23405,22 @@)
$ws.Range("H39").Value = 2838
$ws.Range("I39").Value = 2838
$ws.Range("K39").Value = 2838
$ws.Range("M39").Value = -2447
# row edit (@@ -23898,22 +23901,22 @@)
$ws.Range("H49").Value = 2838
$ws.Range("I49").Value = 2838
$ws.Range("K49").Value = 2838
$ws.Range("M49").Value = -2656
# row edit (@@ -27959,22 +27962,22 @@)
$ws.Range("H132").Value = 64995.9
$ws.Range("I132").Value = 41989.68
$ws.Range("K132").Value = 125969.04
$ws.Range("M132").Value = -123439.04

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# row edit (@@ -30122,25 +30125,25 @@)
$ws.Range("H33").Value = 6328.75
$ws.Range("J33").Value = 25050.5
$ws.Range("L33").Value = 150303
$ws.Range("N33").Value = -150869
# row edit (@@ -30682,25 +30685,25 @@)
$ws.Range("H44").Value = 333666.66
$ws.Range("I44").Value = 500
$ws.Range("J44").Value = 500250
$ws.Range("K44").Value = 1500
$ws.Range("L44").Value = 1500750
$ws.Range("M44").Value = -1102
$ws.Range("N44").Value = -1501546
# row edit (@@ -31934,25 +31937,25 @@)
$ws.Range("H69").Value = 4899.5713
$ws.Range("I69").Value = 4880
$ws.Range("J69").Value = 4902.8335
$ws.Range("K69").Value = 14640
$ws.Range("L69").Value = 14708.5005
$ws.Range("M69").Value = -13829
$ws.Range("N69").Value = -16330.5005
# row edit (@@ -32084,25 +32087,25 @@)
$ws.Range("H72").Value = 4899.5713
$ws.Range("I72").Value = 4880
$ws.Range("J72").Value = 4902.8335
$ws.Range("K72").Value = 43920
$ws.Range("L72").Value = 44125.5015
$ws.Range("M72").Value = -39864
$ws.Range("N72").Value = -52237.5015
# row edit (@@ -32476,22 +32479,22 @@)
$ws.Range("H80").Value = 2500.6667
$ws.Range("I80").Value = 3001
$ws.Range("K80").Value = 9003
$ws.Range("M80").Value = -8067
# row edit (@@ -32629,22 +32632,22 @@)
$ws.Range("H83").Value = 2500.6667
$ws.Range("I83").Value = 3001
$ws.Range("K83").Value = 27009
$ws.Range("M83").Value = -22329
# row edit (@@ -32782,22 +32785,22 @@)
$ws.Range("H86").Value = 456
$ws.Range("I86").Value = 399.5
$ws.Range("K86").Value = 1198.5
$ws.Range("M86").Value = -12.5
# row edit (@@ -32926,22 +32929,22 @@)
$ws.Range("H89").Value = 456
$ws.Range("I89").Value = 399.5
$ws.Range("K89").Value = 3595.5
$ws.Range("M89").Value = 2332.5
# row edit (@@ -34108,25 +34111,25 @@)
$ws.Range("H113").Value = 3415.8262
$ws.Range("I113").Value = 6840
$ws.Range("J113").Value = 2464.6667
$ws.Range("K113").Value = 20520
$ws.Range("L113").Value = 7394.000100000001
$ws.Range("M113").Value = -18350
$ws.Range("N113").Value = -11734.0001
# row edit (@@ -35011,25 +35014,25 @@)
$ws.Range("H131").Value = 16027623
$ws.Range("J131").Value = 23812806
$ws.Range("L131").Value = 71438418
$ws.Range("N131").Value = -71448498

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# row edit (@@ -36684,20 +36687,23 @@)
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 429
# row edit (@@ -37967,22 +37973,19 @@)
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()
# row edit (@@ -38016,22 +38019,22 @@)
$ws.Range("H49").Value = 34949.75
$ws.Range("J49").Value = 34949.75
$ws.Range("L49").Value = 34949.75
$ws.Range("N49").Value = -35317.75

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# row edit (@@ -44500,25 +44503,22 @@)
$ws.Range("H41").Value = 44999
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 44999
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 44999
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -45875
# row edit (@@ -44552,25 +44552,22 @@)
$ws.Range("H42").Value = 8912.5
$ws.Range("I42").Value = 8912.5
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 8912.5
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -8349.5
$ws.Range("N42").ClearContents()
# row edit (@@ -44748,25 +44745,25 @@)
$ws.Range("H46").Value = 6785.2144
$ws.Range("I46").Value = 5363.1816
$ws.Range("J46").Value = 11999.333
$ws.Range("K46").Value = 5363.1816
$ws.Range("L46").Value = 11999.333
$ws.Range("M46").Value = -5175.1816
$ws.Range("N46").Value = -12375.333
# row edit (@@ -44898,25 +44895,22 @@)
$ws.Range("H49").Value = 8912.5
$ws.Range("I49").Value = 8912.5
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 8912.5
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -8765.5
$ws.Range("N49").ClearContents()
# row edit (@@ -47045,22 +47039,22 @@)
$ws.Range("H93").Value = 41669464
$ws.Range("I93").Value = 83335230
$ws.Range("K93").Value = 83335230
$ws.Range("M93").Value = -83333982
# row edit (@@ -48929,25 +48923,25 @@)
$ws.Range("H132").Value = 9063.833000000001
$ws.Range("J132").Value = 5453
$ws.Range("L132").Value = 16359
$ws.Range("N132").Value = -21419

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# row edit (@@ -51247,22 +51241,19 @@)
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()
# row edit (@@ -51893,22 +51884,25 @@)
$ws.Range("H51").Value = 24517.5
$ws.Range("I51").Value = 23070
$ws.Range("J51").Value = 25000
$ws.Range("K51").Value = 23070
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = -22560
$ws.Range("N51").Value = -26020

